$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend columns M (2021) and N (2022) by copying the formatting of column L
# (value + style) and then overwriting with the new figures.
$ws.Range("L4:L12").Copy($ws.Range("M4:M12"))
$ws.Range("L4:L12").Copy($ws.Range("N4:N12"))

# Year headers
$ws.Range("M4").Value = 2021
$ws.Range("N4").Value = 2022

# Data rows
$ws.Range("M5").Value = 5.6
$ws.Range("N5").Value = 6.3

$ws.Range("M6").Value = 0.8
$ws.Range("N6").Value = 0.8

$ws.Range("M7").Value = 1.9
$ws.Range("N7").Value = 2.4

$ws.Range("M8").Value = 0.7
$ws.Range("N8").Value = 0.7

$ws.Range("M9").Value = 0.7
$ws.Range("N9").Value = 0.8

$ws.Range("M10").Value = 0.9
$ws.Range("N10").Value = 1

$ws.Range("M11").Value = 0.3
$ws.Range("N11").Value = 0.2

$ws.Range("M12").Value = 0.2
$ws.Range("N12").Value = 0.4

# New footnote row, matching the style used by the existing footnote cells
$ws.Range("B13").Copy($ws.Range("B14"))
$ws.Range("B14").Value = "По данным лесоустройства 2022 года Лесной службы при Министерстве чрезвычайных ситуаций КР"
$ws.Rows.Item(14).RowHeight = 34.5

# The saved sheet no longer carries a stray cell selection in the view
$ws.Range("A1").Select()
